# Insert a new data row at row 748 (pushing the existing rows 748-789 down
# to 749-790) and populate it with the new record:
#   2026/02/03  火  16  25
#
# Column A holds dates formatted as plain text (e.g. "2026/12/29"), so the
# leading apostrophe forces Excel to keep the value as text instead of
# auto-converting it to a date serial number. Re-applying the "Normal"
# style afterwards clears the transient quote-prefix formatting that the
# apostrophe entry leaves behind, matching the plain, unstyled cells used
# by the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 748..789 down to 749..790, leaving a blank row 748.
$ws.Rows.Item(748).Insert()

# Populate the newly inserted row.
$ws.Cells.Item(748, 1).Value = "'2026/02/03"
$ws.Cells.Item(748, 2).Value = "火"
$ws.Cells.Item(748, 3).Value = 16
$ws.Cells.Item(748, 4).Value = 25

# Clear the quote-prefix style picked up from the text-forced date entry.
$ws.Cells.Item(748, 1).Style = "Normal"
